# Generate Report for Handoff
# Updates the localization-status report: marks the entry as "Ready for
# handoff" (instead of "Handed back: in sync with en-US"), refreshes the
# timestamps, and shrinks columns that previously held the long status
# text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Timestamps
$overview.Range("G2").Value = "2016-09-03 19:13:09"
$zhcn.Range("H2").Value = "2016-09-03 19:13:00"
$dede.Range("H2").Value = "2016-09-03 19:13:09"

# --- Column widths (previously sized for the long status string)
$overview.Columns.Item(5).ColumnWidth = 16.3333333333333
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333
$zhcn.Columns.Item(3).ColumnWidth = 16.3333333333333
$dede.Columns.Item(3).ColumnWidth = 16.3333333333333
